# riska.xlsx — "Add files via upload" data refresh
# Renames the sheet (1) -> (2) and refreshes the numeric/text metrics that
# changed between the two uploaded snapshots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---------------------------------------------
$ws.Name = "repayment_20250912_20250912 (2)"

# --- Helper: write a text-looking value (e.g. "604,330.00") without Excel
# silently re-interpreting it as a formatted number. We flip the cell to
# Text format just long enough to take the literal string, then flip the
# style back to Normal so no stray number-format survives the round trip. --
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- Row 2 : Riska Nurlita --------------------------------------------------
$ws.Range("H2").Value = 210

# --- Row 3 : Romli ----------------------------------------------------------
$ws.Range("D3").Value = 2
Set-TextValue $ws.Range("E3") "604,330.00"
Set-TextValue $ws.Range("G3") "0.40"
$ws.Range("H3").Value = 377
$ws.Range("J3").Value = 1
Set-TextValue $ws.Range("K3") "5.37"
Set-TextValue $ws.Range("L3") "7.14"

# --- Row 4 : Fadilah Damayanti ----------------------------------------------
$ws.Range("H4").Value = 723

# --- Row 5 : Aldi Taufik -----------------------------------------------------
$ws.Range("H5").Value = 899

# --- Row 6 : Axl Wicaksono ----------------------------------------------------
$ws.Range("H6").Value = 465

# --- Row 7 : Annisa Putri Restu ----------------------------------------------
$ws.Range("H7").Value = 473

# --- Row 8 : Debora Retima Sihombing ------------------------------------------
$ws.Range("D8").Value = 4
Set-TextValue $ws.Range("E8") "424,465.00"
Set-TextValue $ws.Range("G8") "0.27"
$ws.Range("H8").Value = 404
$ws.Range("J8").Value = 1
Set-TextValue $ws.Range("K8") "1.43"
Set-TextValue $ws.Range("L8") "6.67"

# --- Row 9 : Azizah Rahmawati -------------------------------------------------
$ws.Range("D9").Value = 1
Set-TextValue $ws.Range("E9") "100,000.00"
Set-TextValue $ws.Range("G9") "0.06"
$ws.Range("H9").Value = 306

# --- Row 10 : Erlangga Hutama -------------------------------------------------
$ws.Range("D10").Value = 1
Set-TextValue $ws.Range("E10") "109,212.00"
Set-TextValue $ws.Range("G10") "0.09"
$ws.Range("H10").Value = 291

# --- Row 11 : Erick Ervan Dewanggga -------------------------------------------
$ws.Range("H11").Value = 272

# --- Row 12 : Ridhoi Berkat Zebua ---------------------------------------------
$ws.Range("H12").Value = 675

# --- Row 14 : Adistira Winditya P ---------------------------------------------
$ws.Range("D14").Value = 4
Set-TextValue $ws.Range("E14") "1,788,670.00"
Set-TextValue $ws.Range("G14") "1.19"
$ws.Range("H14").Value = 203

# --- Row 15 : Yandi Nugraha -----------------------------------------------------
$ws.Range("H15").Value = 1.391

# --- Row 16 : Sucika Wardani -----------------------------------------------------
$ws.Range("D16").Value = 2
Set-TextValue $ws.Range("E16") "371,983.00"
Set-TextValue $ws.Range("G16") "0.26"
$ws.Range("H16").Value = 276
$ws.Range("J16").Value = 1
Set-TextValue $ws.Range("K16") "3.41"
Set-TextValue $ws.Range("L16") "6.67"

# --- Row 17 : Wasti Feronika Sihombing --------------------------------------------
$ws.Range("H17").Value = 305

# --- Row 18 : Nuraini -------------------------------------------------------------
$ws.Range("H18").Value = 158
